$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 2
$ws.Range("F14").Value = -1
$ws.Range("F18").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F28").Value = -1
$ws.Range("F29").Value = 3
$ws.Range("F31").Value = 3
$ws.Range("F33").Value = -4
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = -1
